$d = $word.ActiveDocument

# --- Change 1: collapse the "payors" run/proofErr split into a single run ---
# Paragraph 3 currently holds three runs (plus spell-check proofErr markers)
# around the word "payors"; the final text should be one contiguous run.
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
[void]$r3.MoveEnd(1, -1)
$r3.Delete()
$newText = "Is there a desire to filter merged data based on consumer type: " + `
    [char]8220 + "for example, payors may have different data interest from providers" + `
    [char]8221 + "?"
$r3.InsertAfter($newText)

# --- Change 2: move the _GoBack bookmark from the end of paragraph 5
#     ("Should automatic data validation...") down to a brand-new list
#     paragraph inserted right after paragraph 11
#     ("Are documents expected to be stored...") and before the trailing
#     empty paragraph. ---

# Paragraph 11 is "Are documents expected to be stored..." (still in its
# original position - only the bookmark actually relocates).
$pLast = $d.Paragraphs.Item(11)
$rLast = $pLast.Range
[void]$rLast.Collapse(0)
[void]$rLast.InsertParagraphAfter()

# The freshly inserted paragraph (now index 12) inherits the ListParagraph /
# numPr formatting from paragraph 11 automatically.
$pNew = $d.Paragraphs.Item(12)
$rNew = $pNew.Range
[void]$rNew.MoveEnd(1, -1)
# Give it transient text so we can anchor a bookmark range strictly inside
# this paragraph (an all-empty paragraph's range is ambiguous with the
# following paragraph's start), then strip the placeholder back out.
$rNew.Text = "X"

$rAnchor = $pNew.Range
[void]$rAnchor.MoveEnd(1, -1)

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

[void]$d.Bookmarks.Add("_GoBack", $rAnchor)

$rClear = $pNew.Range
[void]$rClear.MoveEnd(1, -1)
$rClear.Text = ""
